$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 8 (shifts current rows 8-14 down to 9-15)
$ws.Rows.Item(8).Insert()

# Update row 9 (GameStartPing, previously row 8) column E (Datenstruktur)
$ws.Cells.Item(9, 5).Value = "Dict mit data: Repräsentation der Spielerdaten (Objekt der Klasse Player) und players : liste der Spielernamen"

# Update row 7 (UserNameValidationPing) column E (Datenstruktur)
$ws.Cells.Item(7, 5).Value = "Dict mit valid " + [char]0x2208 + " {""True"", ""False""), error " + [char]0x2208 + " {"""", ""doppelt"", ""late""} und players : Liste der Spieler, welche in der Lobby anwesend sind"

# New row 8: NewLobbyPing
$ws.Cells.Item(8, 1).Value = "NewLobbyPing"
$ws.Cells.Item(8, 2).Value = "Server"
$ws.Cells.Item(8, 3).Value = "Client"
$ws.Cells.Item(8, 4).Value = "Sendet dem Client eine Liste der in der Lobby anwesenden Spieler"
$ws.Cells.Item(8, 5).Value = "Liste der Spielernamen"

# Update selection to E8
$ws.Range("E8").Select()
